$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A11 holds a date-like string; force it to stay text instead of being
# auto-converted to a date serial number by entering it via a "@" (text)
# number format, then clearing the format back to the default before
# re-applying the same centered alignment used by the rest of the table.
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2025/11/20"
$ws.Range("A11").ClearFormats()

$ws.Range("B11").Value = "逃离鸭科夫"
$ws.Range("C11").Value = 1196

$ws.Range("A11:C11").HorizontalAlignment = -4108
$ws.Range("A11:C11").VerticalAlignment = -4108
